$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts the string into a
# floating point number and the original textual representation
# (e.g. trailing zero in "74.30") would be lost.
$textRefs = @("D5","D6","D7","D9","D10","D11","D14","D15","D20","D22","D23","D25","D26","D27","D28","D29","D32","D36","D37","D38","D41","D42","D43","D50")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.171.79"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.243.37"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "246.53"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").Value = "74.30"
$ws.Range("E7").Value = "  -3.41%  "
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  -5.57%  "
$ws.Range("D10").Value = "42.21"
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  -3.72%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "14.49"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").Value = "0.850"
$ws.Range("D16").Value = "2.231.66"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "42.059.89"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "0.0₃0983"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "72.17"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").Value = "229.56"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  +41.14%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "11.48"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "169.08"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.11"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("D32").Value = "31.09"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  +9.84%  "
$ws.Range("D36").Value = "4.48"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "0.0312"
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("D38").Value = "13.70"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "62.39"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.203"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "107.13"
$ws.Range("E43").Value = "  -5.37%  "
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").Value = "4.17"
$ws.Range("E50").Value = "  -7.63%  "
$ws.Range("E51").Value = "  +0.42%  "
